# Update "想去人数" (want-to-go count) figures in column F across all four
# sheets to match the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2773
$ws.Range("F5").Value = 948
$ws.Range("F7").Value = 2442
$ws.Range("F8").Value = 1879
$ws.Range("F11").Value = 2533
$ws.Range("F16").Value = 136
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 9458
$ws.Range("F21").Value = 7422
$ws.Range("F22").Value = 11963
$ws.Range("F27").Value = 575
$ws.Range("F28").Value = 2690
$ws.Range("F30").Value = 212
$ws.Range("F31").Value = 2667
$ws.Range("F32").Value = 967
$ws.Range("F37").Value = 1059
$ws.Range("F40").Value = 57
$ws.Range("F41").Value = 560

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 7

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 182

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2773
$ws.Range("F8").Value = 948
$ws.Range("F11").Value = 2442
$ws.Range("F13").Value = 1879
$ws.Range("F15").Value = 2533
$ws.Range("F20").Value = 136
$ws.Range("F21").Value = 124
$ws.Range("F22").Value = 9458
$ws.Range("F25").Value = 7422
$ws.Range("F26").Value = 11963
$ws.Range("F32").Value = 575
$ws.Range("F34").Value = 2690
$ws.Range("F37").Value = 212
$ws.Range("F45").Value = 560
